$d = $word.ActiveDocument

# Locate the "John Choy" run in the title paragraph so we can compute the
# split point between "John " and "Choy" without hard-coding offsets.
$findRange = $d.Content
$found = $findRange.Find.Execute("John Choy", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Bookmark span: from the very start of the document/paragraph through the
# end of "John " (i.e. right before "Choy"). Adding a bookmark whose range
# ends mid-run causes Word to split that run, matching the target XML.
$splitPoint = $findRange.Start + 5  # "John " is 5 characters

$bmRange1 = $d.Range(0, $splitPoint)
$d.Bookmarks.Add("OLE_LINK1", $bmRange1)

$bmRange2 = $d.Range(0, $splitPoint)
$d.Bookmarks.Add("OLE_LINK2", $bmRange2)
